# Update excess mortality analysis to week 35
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update revised "Waargenomen" (observed) values in column G ---
# These corrections ripple into column I ("Oversterfte") via the
# existing shared formula G-H, which Excel will recalc automatically.
$ws.Range("G4").Value  = 3613
$ws.Range("G7").Value  = 4978
$ws.Range("G15").Value = 2682
$ws.Range("G21").Value = 2523
$ws.Range("G22").Value = 2669
$ws.Range("G23").Value = 2651
$ws.Range("G24").Value = 2628
$ws.Range("G25").Value = 3197
$ws.Range("G26").Value = 2821

# --- Move the totals row from row 28 down to row 29, freeing row 27 ---
$ws.Range("F28").Value = $null
$ws.Range("G28").Value = $null
$ws.Range("H28").Value = $null
$ws.Range("I28").Value = $null

# --- Add new data row 27 for week 35 ---
$ws.Range("F27").Value = 35
$ws.Range("G27").Value = 2689
$ws.Range("H27").Value = 2822
$ws.Range("I27").Formula = "=G27-H27"

# --- Recreate the totals row at row 29, now summing through row 27 ---
$ws.Range("F29").Value = "Som week 11 tot en met 19"
$ws.Range("G29").Formula = "=SUM(G3:G27)"
$ws.Range("H29").Formula = "=SUM(H3:H27)"
$ws.Range("I29").Formula = "=SUM(I3:I27)"

$ws.Range("G29:I29").Style = "Comma"

# --- Update the view: scroll so row 17 is at top, select the new total row ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("G29:I29").Select()

$wb.Application.Calculate()
